$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OPEX")

$ws.Range("B3").Formula = "=1.5*19283.0421189704"
$ws.Range("B4").Formula = "=1.5*43709.8236922189"
$ws.Range("B5").Formula = "=1.5*44248.4220525376"
$ws.Range("B6").Formula = "=1.5*10540.7765197658"
$ws.Range("B7").Formula = "=1.5*11938.9186890398"
$ws.Range("B8").Formula = "=1.5*50983.7161218294"
$ws.Range("B9").Formula = "=1.5*47726.863400003"
$ws.Range("B10").Formula = "=1.5*48220.2891692096"
$ws.Range("B11").Formula = "=1.5*37682.86771352"
$ws.Range("B12").Formula = "=1.5*44827.8721002112"
$ws.Range("B13").Formula = "=1.5*12139.8149301997"
$ws.Range("B14").Formula = "=1.5*47419.5747604423"
$ws.Range("B15").Formula = "=1.5*50088.4769205996"
$ws.Range("B2").Formula = "=AVERAGE(B3:B15)"

$ws.Activate()
$ws.Range("B2").Select()
